$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1985
$ws.Range("B4").Value = 1941
$ws.Range("B1:B4").NumberFormat = "General"
$ws.Range("B4").Select()
